$wb = $excel.ActiveWorkbook

# --- Sheet "general": update objValue, runtime, Z1-Z4 ---
$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 38.17098846897687
$ws.Range("B4").Value = 0.01199984550476074
$ws.Range("B6").Value = 38.17098846897688
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0

# --- Sheet "x": update permutation column B ---
$ws = $wb.Worksheets.Item("x")
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 9
$ws.Range("B4").Value = 10
$ws.Range("B5").Value = 12
$ws.Range("B6").Value = 6
$ws.Range("B8").Value = 3
$ws.Range("B9").Value = 2
$ws.Range("B10").Value = 11
$ws.Range("B11").Value = 1
$ws.Range("B13").Value = 8
$ws.Range("B14").Value = 7

# --- Sheet "TBar": update column B ---
$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value = 32.61192465059682
$ws.Range("B4").Value = 32.71621059566448
$ws.Range("B5").Value = 39.77511225757775
$ws.Range("B6").Value = 34.76592070603971
$ws.Range("B8").Value = 37.06506101847738
$ws.Range("B9").Value = 32.01159140980468
$ws.Range("B10").Value = 32.45367071955468
$ws.Range("B11").Value = 30
$ws.Range("B12").Value = 30
$ws.Range("B13").Value = 39.73013137402148
$ws.Range("B14").Value = 37.94859027624736
$ws.Range("B15").Value = 42.09384035720478

# --- Sheet "y": remove data rows (keep header only) ---
$ws = $wb.Worksheets.Item("y")
$ws.Range("A2:D6").ClearContents() | Out-Null

# --- Sheet "Q": update column C ---
$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value = 250.970000000001
$ws.Range("C8").Value = 260.990000000001
$ws.Range("C9").Value = 252.9750000000009
$ws.Range("C10").Value = 269.580000000001
$ws.Range("C11").Value = 250.575000000001
$ws.Range("C12").Value = 151.3550000000009
$ws.Range("C13").Value = 157
$ws.Range("C14").Value = 157.125000000001
$ws.Range("C15").Value = 153.080000000001
$ws.Range("C16").Value = 160.235000000001
$ws.Range("C17").Value = 272.9599999999988
$ws.Range("C18").Value = 283.2849999999988
$ws.Range("C19").Value = 256.7699999999988
$ws.Range("C20").Value = 275.9449999999989
$ws.Range("C21").Value = 263.9099999999989
$ws.Range("C22").Value = 226.0399999999994
$ws.Range("C23").Value = 247.1799999999994
$ws.Range("C24").Value = 221.8549999999994
$ws.Range("C25").Value = 238.4549999999994
$ws.Range("C26").Value = 224.4749999999995
$ws.Range("C32").Value = 258.7350000000008
$ws.Range("C33").Value = 269.2400000000008
$ws.Range("C34").Value = 250.9150000000008
$ws.Range("C35").Value = 261.9150000000008
$ws.Range("C36").Value = 255.0150000000008
$ws.Range("C37").Value = 141.0250000000001
$ws.Range("C38").Value = 143.4
$ws.Range("C39").Value = 139.7050000000001
$ws.Range("C40").Value = 150.4250000000002
$ws.Range("C41").Value = 134.7700000000002
$ws.Range("C42").Value = 153.4099999999999
$ws.Range("C43").Value = 167.1249999999999
$ws.Range("C44").Value = 139.5349999999999
$ws.Range("C45").Value = 154.5
$ws.Range("C46").Value = 143.6599999999999
$ws.Range("C47").Value = 85.48500000000051
$ws.Range("C48").Value = 87.9650000000005
$ws.Range("C49").Value = 79.71500000000052
$ws.Range("C50").Value = 90.7300000000005
$ws.Range("C51").Value = 84.73000000000052
$ws.Range("C52").Value = 57.95
$ws.Range("C53").Value = 58.67999999999927
$ws.Range("C54").Value = 61.72999999999927
$ws.Range("C55").Value = 60.65499999999928
$ws.Range("C56").Value = 52.91499999999927
$ws.Range("C57").Value = 258.7350000000008
$ws.Range("C58").Value = 269.2400000000008
$ws.Range("C59").Value = 250.9150000000008
$ws.Range("C60").Value = 261.9150000000008
$ws.Range("C61").Value = 255.0150000000008
$ws.Range("C62").Value = 250.970000000001
$ws.Range("C63").Value = 260.990000000001
$ws.Range("C64").Value = 252.9750000000009
$ws.Range("C65").Value = 269.580000000001
$ws.Range("C66").Value = 250.575000000001
$ws.Range("C67").Value = 272.9599999999988
$ws.Range("C68").Value = 283.2849999999988
$ws.Range("C69").Value = 256.7699999999988
$ws.Range("C70").Value = 275.9449999999989
$ws.Range("C71").Value = 263.9099999999989

# --- Sheet "R": zero out rows 7-11 (j=12 group) ---
$ws = $wb.Worksheets.Item("R")
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 0

# --- Sheet "rho": remove data rows (keep header only) ---
$ws = $wb.Worksheets.Item("rho")
$ws.Range("A2:C6").ClearContents() | Out-Null

# --- Sheet "alpha": remove data rows (keep header only) ---
$ws = $wb.Worksheets.Item("alpha")
$ws.Range("A2:C6").ClearContents() | Out-Null

Write-Host "edit complete"
